$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# Overview sheet: E2, F2 status; G2 timestamp
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
$wsOverview.Range("G2").Value = "2016-08-26 22:57:45"

# zh-cn sheet: C2 status; H2 timestamp
$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsZhCn.Range("H2").Value = "2016-08-26 22:57:41"

# de-de sheet: C2 status; H2 timestamp
$wsDeDe.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("H2").Value = "2016-08-26 22:57:45"

# Column width adjustments (status columns widened to fit "Ready for handoff")
# NB: ColumnWidth is quantized to whole pixels (1/6-character-unit steps) by
# this host, so 16.333333333333332 is the input that lands on the bucket
# nearest the canonical target width (~17.22 chars -> 17.166666666666668).
$wsOverview.Range("E1:F1").ColumnWidth = 16.333333333333332
$wsZhCn.Range("C1").ColumnWidth = 16.333333333333332
$wsDeDe.Range("C1").ColumnWidth = 16.333333333333332
